$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $c = $ws.Range($cell)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextCell "D2" "328.89"
Set-TextCell "E2" "1.64%"
Set-TextCell "D3" "41.05"
Set-TextCell "E3" "3.27%"
Set-TextCell "D4" "5.591"
Set-TextCell "E4" "-4.93%"
Set-TextCell "D5" "0.08174"
Set-TextCell "E5" "1.73%"
Set-TextCell "D6" "2.034"
Set-TextCell "E6" "5.62%"
Set-TextCell "D7" "8.738"
Set-TextCell "E7" "0.77%"
Set-TextCell "D8" "4.534"
Set-TextCell "E8" "-1.18%"
Set-TextCell "E9" "1.58%"
Set-TextCell "D10" "0.9183"
Set-TextCell "E10" "-1.61%"
Set-TextCell "D11" "0.1262"
Set-TextCell "E11" "-0.83%"
Set-TextCell "D12" "0.1959"
Set-TextCell "E12" "-0.50%"
Set-TextCell "D13" "0.09393"
Set-TextCell "E13" "1.95%"
Set-TextCell "D14" "0.03736"
Set-TextCell "E14" "5.09%"
Set-TextCell "D15" "0.1057"
Set-TextCell "E15" "1.06%"
Set-TextCell "D16" "0.001299"
Set-TextCell "D17" "0.006302"
Set-TextCell "E17" "2.46%"
Set-TextCell "D18" "3.437"
Set-TextCell "E18" "2.76%"
Set-TextCell "E19" "-2.27%"
Set-TextCell "D20" "8.316"
Set-TextCell "E20" "-4.81%"
Set-TextCell "D21" "0.1393"
Set-TextCell "E21" "-1.77%"
Set-TextCell "E22" "2.82%"
Set-TextCell "D23" "0.04423"
Set-TextCell "E23" "0.25%"
Set-TextCell "D24" "0.001262"
Set-TextCell "E24" "0.12%"
Set-TextCell "D25" "0.004300"
Set-TextCell "E25" "-3.14%"
Set-TextCell "D39" "0.02757"
Set-TextCell "E39" "12.77%"
Set-TextCell "D40" "0.05405"
Set-TextCell "E40" "2.94%"
Set-TextCell "D41" "0.007671"
Set-TextCell "E41" "3.33%"
Set-TextCell "D42" "0.1414"
Set-TextCell "E42" "0.56%"
Set-TextCell "D43" "0.008998"
Set-TextCell "E43" "-5.82%"
Set-TextCell "D44" "0.002124"
Set-TextCell "E44" "0.33%"
Set-TextCell "D45" "0.01133"
Set-TextCell "E45" "13.76%"
Set-TextCell "D46" "0.00006909"
Set-TextCell "E46" "2.60%"
Set-TextCell "D47" "0.00000000752"
Set-TextCell "B48" "CoinbaseStockToken"
Set-TextCell "C48" "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
Set-TextCell "D48" "0.002284"
Set-TextCell "E48" "60.55%"
Set-TextCell "B49" "BOLO"
Set-TextCell "C49" "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
Set-TextCell "D49" "0.003587"
Set-TextCell "E49" "19.48%"
